$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.237.88'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '1.905.73'
$ws.Range("E3").Value = '  +1.80%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  -0.03%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.82'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("E6").Value = '  +0.07%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4628'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +0.22%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3955'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  +1.65%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.75'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +1.09%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07968'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +1.17%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.000'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +2.54%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.30'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").Value = '1.934.64'
$ws.Range("E13").Value = '  +4.41%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.133'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +1.79%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.761'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +0.99%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06953'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -0.03%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.77'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = '29.268.09'
$ws.Range("E22").Value = '  +1.96%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.361'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +1.71%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '2.155.47'
$ws.Range("E25").Value = '  +4.41%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.048'
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.80'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +2.72%  '
$ws.Range("E28").Value = '  +1.16%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.919'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +0.81%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.000'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("E31").Value = '  -0.05%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09410'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.85%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9252'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("E35").Value = '  +1.22%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.260'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -1.90%  '
$ws.Range("E37").Value = '  +0.63%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.168'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +1.51%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02109'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +1.56%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.000'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +4.17%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5759'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  +2.34%  '
$ws.Range("E42").Value = '  +1.37%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.988'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +2.12%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.03'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +2.80%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5431'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +2.73%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.226'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +2.95%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07105'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -1.63%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.879'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +2.14%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.577'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +6.85%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.00'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -0.72%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.060'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -5.85%  '
